# Apply "adjust for difference instance" edits to schedule worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..21 (columns B, C, D, E, F). A, G, H are unchanged.
$data = @(
    @{Row=2;  B=1;  C=0.02222222222222222; D=0.04444444444444445; E="C19"; F=2},
    @{Row=3;  B=2;  C=0.01041666666666667; D=0.03472222222222222; E="F56"; F=1},
    @{Row=4;  B=3;  C=0.01041666666666667; D=0.03402777777777777; E="A9";  F=2},
    @{Row=5;  B=4;  C=0.01875;             D=0.04305555555555556; E="F54"; F=3},
    @{Row=6;  B=5;  C=0.01388888888888889; D=0.03541666666666667; E="A7";  F=1},
    @{Row=7;  B=6;  C=0.02083333333333333; D=0.04791666666666667; E="D49"; F=2},
    @{Row=8;  B=7;  C=0.06736111111111111; D=0.08958333333333333; E="B1";  F=2},
    @{Row=9;  B=8;  C=0.06041666666666667; D=0.08749999999999999; E="D42"; F=1},
    @{Row=10; B=9;  C=0.07152777777777777; D=0.09305555555555556; E="A17"; F=1},
    @{Row=11; B=10; C=0.04930555555555555; D=0.07291666666666667; E="C2";  F=2},
    @{Row=12; B=11; C=0.06527777777777778; D=0.09027777777777778; E="B2";  F=1},
    @{Row=13; B=12; C=0.06458333333333334; D=0.09166666666666666; E="E11"; F=1},
    @{Row=14; B=13; C=0.0625;              D=0.08749999999999999; E="F31"; F=1},
    @{Row=15; B=14; C=0.05486111111111111; D=0.0763888888888889;  E="F33"; F=1},
    @{Row=16; B=15; C=0.07777777777777778; D=0.1020833333333333;  E="C14"; F=3},
    @{Row=17; B=16; C=0.05486111111111111; D=0.07847222222222222; E="C17"; F=3},
    @{Row=18; B=17; C=0.04236111111111111; D=0.06458333333333334; E="F50"; F=2},
    @{Row=19; B=18; C=0.04166666666666666; D=0.0625;               E="E23"; F=3},
    @{Row=20; B=19; C=0.05069444444444444; D=0.07708333333333334; E="C16"; F=2},
    @{Row=21; B=20; C=0.05694444444444444; D=0.07777777777777778; E="E7";  F=3}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}
